$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2").Value = "DemoSalesManager"
$ws.Range("A4").Select()
